$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stats for the latest month (row 23) to reflect the 2025-10 refresh
$ws.Range("B23").Value = 6314
$ws.Range("C23").Value = 1004
$ws.Range("D23").Value = 5912688
$ws.Range("E23").Value = 936.4409249287298
$ws.Range("F23").Value = 8.339052848318463
$ws.Range("G23").Value = 4.474505723205002
$ws.Range("H23").Value = 26.69376952173739
